# "17th April 1st update" — add the 17/04/2020 data point to the national
# timeseries sheet and revise the previous day's (16/04/2020) running total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34 is 16/04/2020 — its case count is revised from 660 to 1062.
$ws.Range("B34").Value = 1062

# Insert a new row 36 (after 17/03/2020 on row 35, before what is currently
# row 36 / 18/03/2020) for the new date 17/04/2020, pushing every row below
# down by one.
$ws.Rows.Item(36).Insert()
$ws.Range("A36").Value = "17/04/2020"
$ws.Range("B36").Value = 231
